# TOR-012 Update data sheets
# Updates existing row 3 figures and appends two new roster rows (4 & 5).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 3 values ---
$ws.Range("C3").Value = 1
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 12

# --- Add new row 4 ---
$ws.Range("A4").Value = 5
$ws.Range("B4").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 3
$ws.Range("G4").Value = 8
$ws.Range("A4:G4").Locked = $false

# --- Add new row 5 ---
$ws.Range("A5").Value = 6
$ws.Range("B5").Value = 2
$ws.Range("C5").Value = 2
$ws.Range("D5").Value = 2
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 3
$ws.Range("G5").Value = 6
$ws.Range("B5:G5").Locked = $false

# --- Update selection to match the final cursor position ---
$ws.Range("I7").Select()
